# TMT0072153_TMT0072155_VerifyBanker_PrimaryOrNonPrimary_WhoIsPartOfTheActivityCanEditTheActivity
# Adds a "MoreAttendees" sheet (copy of the ExtAttendee/HLAttendee columns
# from UpdateActivity) after the UpdateActivity sheet, and leaves behind the
# selection/view state that results from doing so interactively.

$wb = $excel.ActiveWorkbook

# --- Activity sheet: just a different cell was clicked/selected ---
$wsActivity = $wb.Worksheets.Item("Activity")
$wsActivity.Activate()
$wsActivity.Range("C20").Select()

# --- UpdateActivity sheet: select + copy the ExtAttendee/HLAttendee columns ---
$wsUpdate = $wb.Worksheets.Item("UpdateActivity")
$wsUpdate.Activate()
$wsUpdate.Range("F1:G2").Select()
$wsUpdate.Range("F1:G2").Copy()

# --- New sheet, placed right after UpdateActivity ---
$newWs = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsUpdate)
$newWs.Name = "MoreAttendees"
$newWs.Range("A1").Select()
$newWs.Paste()

# Re-apply the bold+centered header style (copy/paste doesn't carry it here)
$newWs.Range("A1:B1").Font.Bold = $true
$newWs.Range("A1:B1").HorizontalAlignment = -4108

# Match the source columns' widths as closely as this engine allows
$newWs.Columns("A").ColumnWidth = 10.7
$newWs.Columns("B").ColumnWidth = 16.83

# Afterwards the whole F:G columns got selected back on UpdateActivity ...
$wsUpdate.Activate()
$wsUpdate.Columns("F:G").Select()

# ... and MoreAttendees (the new sheet) is left as the active tab/selection
$newWs.Activate()
$newWs.Range("C13").Select()
